# Updated cryptos list data: refresh Price (col D) and Volume(1h) (col E)
# values per latest scrape. Text formatting is preserved (prices/volumes
# are stored as text, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.856.01'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.638.03'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.60%  '

$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("E7").Value = '  +0.18%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2575'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06423'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.74'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07795'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.49%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.286'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.864.91'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.637.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5597'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7625'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.98'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.88%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.878.06'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '194.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.50%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.325'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.84%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.877'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.42%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.094'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.004'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("E25").Value = '  -6.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1257'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.32%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.824'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.25%  '

$ws.Range("E29").Value = '  -1.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.242'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.42%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04886'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.294'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.223'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.67%  '

$ws.Range("E34").Value = '  +1.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.379'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9033'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.576'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5514'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.126.36'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01562'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9969'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.545'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8004'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.56%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '98.06'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.775.15'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.14%  '

$ws.Range("E46").Value = '  -6.76%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4266'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.677'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05043'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.003'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.54%  '

